# Auto-generated: apply Leviathan_Profits price/profit recalculation updates
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 220.1
$ws.Range("I9").Value = 327.54544
$ws.Range("K9").Value = 327.54544
$ws.Range("M9").Value = -158.54544

$ws.Range("H28").Value = 133
$ws.Range("I28").Value = 133
$ws.Range("K28").Value = 133
$ws.Range("M28").Value = 352

$ws.Range("H38").Value = 401.69232
$ws.Range("I38").Value = 401.69232
$ws.Range("K38").Value = 1205.07696
$ws.Range("M38").Value = -833.0769599999999

$ws.Range("H93").Value = 30000
$ws.Range("J93").Value = 30000
$ws.Range("L93").Value = 30000
$ws.Range("N93").Value = -34992

$ws.Range("H98").Value = 1504.0333
$ws.Range("I98").Value = 1525.5172
$ws.Range("K98").Value = 1525.5172
$ws.Range("M98").Value = -27.5172

$ws.Range("H122").Value = 1504.0333
$ws.Range("I122").Value = 1525.5172
$ws.Range("K122").Value = 4576.5516
$ws.Range("M122").Value = -2126.5516

$ws.Range("H132").Value = 1843.4054
$ws.Range("I132").Value = 1903.1177
$ws.Range("K132").Value = 5709.3531
$ws.Range("M132").Value = -3179.3531

$ws.Range("H137").Value = 4662.9434
$ws.Range("I137").Value = 3862.3057
$ws.Range("J137").Value = 6358.4116
$ws.Range("K137").Value = 11586.9171
$ws.Range("L137").Value = 19075.2348
$ws.Range("M137").Value = -9036.917099999999
$ws.Range("N137").Value = -24175.2348

$ws.Range("H138").Value = 1602.2727
$ws.Range("J138").Value = 3489.077
$ws.Range("L138").Value = 10467.231
$ws.Range("N138").Value = -20747.231

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 24205.438
$ws.Range("I32").Value = 4475.582
$ws.Range("J32").Value = 144776.78
$ws.Range("K32").Value = 4475.582
$ws.Range("L32").Value = 144776.78
$ws.Range("M32").Value = -4188.582
$ws.Range("N32").Value = -145350.78

$ws.Range("H61").Value = 994.1
$ws.Range("I61").Value = 882.3333
$ws.Range("K61").Value = 882.3333
$ws.Range("M61").Value = -670.3333

$ws.Range("H63").Value = 1811.9
$ws.Range("I63").Value = 1811.9
$ws.Range("K63").Value = 1811.9
$ws.Range("M63").Value = -1125.9

$ws.Range("H66").Value = 1811.9
$ws.Range("I66").Value = 1811.9
$ws.Range("K66").Value = 9059.5
$ws.Range("M66").Value = -5627.5

$ws.Range("H136").Value = 994.1
$ws.Range("I136").Value = 882.3333
$ws.Range("K136").Value = 2646.9999
$ws.Range("M136").Value = -96.9998999999998

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H27").Value = 30542
$ws.Range("J27").Value = 30542
$ws.Range("L27").Value = 30542
$ws.Range("N27").Value = -30926

$ws.Range("H105").Value = 3985.2942
$ws.Range("I105").Value = 4330
$ws.Range("K105").Value = 4330
$ws.Range("M105").Value = -2583

$ws.Range("H132").Value = 76330
$ws.Range("J132").Value = 76330
$ws.Range("L132").Value = 76330
$ws.Range("N132").Value = -86450

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 45265.695
$ws.Range("I31").Value = 53994.633
$ws.Range("J31").Value = 3803.25
$ws.Range("K31").Value = 53994.633
$ws.Range("L31").Value = 3803.25
$ws.Range("M31").Value = -53699.633
$ws.Range("N31").Value = -4393.25

$ws.Range("H34").Value = 45265.695
$ws.Range("I34").Value = 53994.633
$ws.Range("J34").Value = 3803.25
$ws.Range("K34").Value = 53994.633
$ws.Range("L34").Value = 3803.25
$ws.Range("M34").Value = -53792.633
$ws.Range("N34").Value = -4207.25

$ws.Range("H58").Value = 1594.9678
$ws.Range("I58").Value = 1497.6207
$ws.Range("K58").Value = 1497.6207
$ws.Range("M58").Value = -1294.6207

$ws.Range("H122").Value = 5027.8887
$ws.Range("I122").Value = 3700.1538
$ws.Range("K122").Value = 11100.4614
$ws.Range("M122").Value = -8650.4614

$ws.Range("H136").Value = 1594.9678
$ws.Range("I136").Value = 1497.6207
$ws.Range("K136").Value = 4492.8621
$ws.Range("M136").Value = -1942.8621

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 108.35294
$ws.Range("J12").Value = 176.6
$ws.Range("L12").Value = 529.8
$ws.Range("N12").Value = -875.8

$ws.Range("H68").Value = 983.44446
$ws.Range("J68").Value = 684
$ws.Range("L68").Value = 2052
$ws.Range("N68").Value = -3674

$ws.Range("H71").Value = 983.44446
$ws.Range("J71").Value = 684
$ws.Range("L71").Value = 6156
$ws.Range("N71").Value = -14268

$ws.Range("H96").Value = 0
$ws.Range("I96").Value = 0
$ws.Range("J96").Value = 0
$ws.Range("K96").Value = 0
$ws.Range("N96").ClearContents()
$ws.Range("M96").ClearContents()

$ws.Range("H104").Value = 251113.25
$ws.Range("I104").Value = 999999
$ws.Range("K104").Value = 2999997
$ws.Range("M104").Value = -2997376

$ws.Range("H115").Value = 0
$ws.Range("I115").Value = 0
$ws.Range("J115").Value = 0
$ws.Range("K115").Value = 0
$ws.Range("N115").ClearContents()
$ws.Range("M115").ClearContents()

$ws.Range("H117").Value = 4639.6816
$ws.Range("I117").Value = 2202
$ws.Range("K117").Value = 6606
$ws.Range("M117").Value = -3164

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H36").Value = 16229
$ws.Range("J36").Value = 29998
$ws.Range("L36").Value = 29998
$ws.Range("N36").Value = -30968

$ws.Range("H70").Value = 5152.1
$ws.Range("J70").Value = 5066.6665
$ws.Range("L70").Value = 5066.6665
$ws.Range("N70").Value = -5606.6665

$ws.Range("H73").Value = 5152.1
$ws.Range("J73").Value = 5066.6665
$ws.Range("L73").Value = 5066.6665
$ws.Range("N73").Value = -6938.6665

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 1919.6
$ws.Range("I68").Value = 1513.7142
$ws.Range("J68").Value = 2274.75
$ws.Range("K68").Value = 1513.7142
$ws.Range("L68").Value = 2274.75
$ws.Range("M68").Value = -764.7141999999999
$ws.Range("N68").Value = -3772.75

$ws.Range("H71").Value = 1919.6
$ws.Range("I71").Value = 1513.7142
$ws.Range("J71").Value = 2274.75
$ws.Range("K71").Value = 7568.571
$ws.Range("L71").Value = 11373.75
$ws.Range("M71").Value = -3824.571
$ws.Range("N71").Value = -18861.75

$ws.Range("H81").Value = 20164
$ws.Range("I81").Value = 20164
$ws.Range("K81").Value = 20164
$ws.Range("M81").Value = -19166

$ws.Range("H84").Value = 20164
$ws.Range("I84").Value = 20164
$ws.Range("K84").Value = 60492
$ws.Range("M84").Value = -55500

$ws.Range("H132").Value = 3035.2126
$ws.Range("I132").Value = 2651.0286
$ws.Range("J132").Value = 4155.75
$ws.Range("K132").Value = 7953.085800000001
$ws.Range("L132").Value = 12467.25
$ws.Range("M132").Value = -5423.085800000001
$ws.Range("N132").Value = -17527.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 40373
$ws.Range("J54").Value = 46164.332
$ws.Range("L54").Value = 46164.332
$ws.Range("N54").Value = -47204.332

$ws.Range("H81").Value = 6040.857
$ws.Range("I81").Value = 5056.0415
$ws.Range("J81").Value = 11949.75
$ws.Range("K81").Value = 10112.083
$ws.Range("L81").Value = 23899.5
$ws.Range("M81").Value = -9051.083000000001
$ws.Range("N81").Value = -26021.5

$ws.Range("H84").Value = 6040.857
$ws.Range("I84").Value = 5056.0415
$ws.Range("J84").Value = 11949.75
$ws.Range("K84").Value = 50560.415
$ws.Range("L84").Value = 119497.5
$ws.Range("M84").Value = -45256.415
$ws.Range("N84").Value = -130105.5

$ws.Range("H113").Value = 620.8
$ws.Range("I113").Value = 467.55554
$ws.Range("K113").Value = 1402.66662
$ws.Range("M113").Value = 767.33338

$ws.Range("H122").Value = 2244.8333
$ws.Range("J122").Value = 2489.6667
$ws.Range("L122").Value = 7469.000100000001
$ws.Range("N122").Value = -12369.0001

$ws.Range("H126").Value = 7416.4375
$ws.Range("I126").Value = 2205.4614
$ws.Range("K126").Value = 6616.3842
$ws.Range("M126").Value = -4146.3842

$ws.Range("H132").Value = 1327309.4
$ws.Range("I132").Value = 5445.893
$ws.Range("K132").Value = 16337.679
$ws.Range("M132").Value = -13807.679
